$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 70; existing rows 70..174 shift down to 71..175
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new data point
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44571
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 100112010
$ws.Cells.Item(70, 7).Value = "Achicoria"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 125
$ws.Cells.Item(70, 11).Value = 6000
$ws.Cells.Item(70, 12).Value = 6500
$ws.Cells.Item(70, 13).Value = 6240
$ws.Cells.Item(70, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(70, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(70, 16).Value = 390
$ws.Cells.Item(70, 17).Value = 16
$ws.Cells.Item(70, 18).Value = "Hortaliza"
